# Update data files - Bot run at 2026-02-13 04:57:55 UTC
# Applies updated API usage counters/percentages to rows 2 and 14 of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (llama-3.1-8b-instant): Current_Ct_Day, Current_Pct_Ct, Current_Ct_Tokens, Current_Pct_Tokens
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 0.0001388888888888889
$ws.Range("K2").Value = 937
$ws.Range("L2").Value = 0.001874

# Row 14 (qwen/qwen3-32b): Current_Ct_Day, Current_Pct_Ct, Current_Ct_Tokens, Current_Pct_Tokens
$ws.Range("I14").Value = 1
$ws.Range("J14").Value = 0.001
$ws.Range("K14").Value = 483
$ws.Range("L14").Value = 0.000966
